# Add BJT Hand calculations into column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.047
$ws.Range("B4").Value = 2.214
$ws.Range("B5").Value = 1.514
$ws.Range("B6").Formula = "=B3-B5"
$ws.Range("B7").Formula = "=186*10^-6"
$ws.Range("B8").Value = 0.0291
$ws.Range("B9").Value = 0.0293
$ws.Range("B13").Value = 282.76
$ws.Range("B17").Value = 157.017
$ws.Range("B18").Value = 318.35

$ws.Range("D19").Select() | Out-Null
